$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final Cluster / Active Cases data (sorted order matching the shared-string table),
# rows 2..63 covering A2:B63. Existing rows are overwritten in place and new rows
# are appended past the old 56-row extent.

$ws.Cells.Item(2, 1).Value = "3323 Villa Maria Catholic Homes St Bernadette's Aged Care Sunshine North"
$ws.Cells.Item(2, 2).Value = 13
$ws.Cells.Item(3, 1).Value = "3364 Assisi Centre Aged Care Rosanna"
$ws.Cells.Item(3, 2).Value = 11
$ws.Cells.Item(4, 1).Value = "3376 Royal Freemasons Coppin Centre Melbourne"
$ws.Cells.Item(4, 2).Value = 33
$ws.Cells.Item(5, 1).Value = "3653 Fronditha Thalpori St Albans Aged Care"
$ws.Cells.Item(5, 2).Value = 19
$ws.Cells.Item(6, 1).Value = "3825 TLC Forest Lodge Residential Aged Care Frankston North"
$ws.Cells.Item(6, 2).Value = 16
$ws.Cells.Item(7, 1).Value = "4167 Royal Freemasons Centennial Lodge Wantirna South"
$ws.Cells.Item(7, 2).Value = 22
$ws.Cells.Item(8, 1).Value = "44226 Boneo Primary School Boneo"
$ws.Cells.Item(8, 2).Value = 11
$ws.Cells.Item(9, 1).Value = "44321 Maiden Gully Primary School Maiden Gully"
$ws.Cells.Item(9, 2).Value = 10
$ws.Cells.Item(10, 1).Value = "44852 Dandenong South Primary School Dandenong"
$ws.Cells.Item(10, 2).Value = 16
$ws.Cells.Item(11, 1).Value = "45034 River Gum Primary School Hampton Park"
$ws.Cells.Item(11, 2).Value = 10
$ws.Cells.Item(12, 1).Value = "45573 Narre Warren South P-12 College Narre Warren South"
$ws.Cells.Item(12, 2).Value = 18
$ws.Cells.Item(13, 1).Value = "45585 Mount Ridley College Craigieburn"
$ws.Cells.Item(13, 2).Value = 10
$ws.Cells.Item(14, 1).Value = "45695 Sacred Heart Primary School Yarrawonga"
$ws.Cells.Item(14, 2).Value = 47
$ws.Cells.Item(15, 1).Value = "4574 Village Glen Aged Care Residences Mornington"
$ws.Cells.Item(15, 2).Value = 10
$ws.Cells.Item(16, 1).Value = "45804 St Therese's School Essendon"
$ws.Cells.Item(16, 2).Value = 10
$ws.Cells.Item(17, 1).Value = "46322 Minaret College Officer Campus Officer"
$ws.Cells.Item(17, 2).Value = 41
$ws.Cells.Item(18, 1).Value = "46390 Al Siraat College Epping"
$ws.Cells.Item(18, 2).Value = 34
$ws.Cells.Item(19, 1).Value = "50516 Ilim College Glenroy Campus Hadfield"
$ws.Cells.Item(19, 2).Value = 12
$ws.Cells.Item(20, 1).Value = "50567 Alamanda K9 College Point Cook"
$ws.Cells.Item(20, 2).Value = 11
$ws.Cells.Item(21, 1).Value = "51478 Wyndham Vale Primary School Wyndham Vale"
$ws.Cells.Item(21, 2).Value = 10
$ws.Cells.Item(22, 1).Value = "52380 Al Iman College Melton South"
$ws.Cells.Item(22, 2).Value = 20
$ws.Cells.Item(23, 1).Value = "52912 Edgars Creek Primary School Wollert"
$ws.Cells.Item(23, 2).Value = 12
$ws.Cells.Item(24, 1).Value = "52985 Minaret College Springvale"
$ws.Cells.Item(24, 2).Value = 17
$ws.Cells.Item(25, 1).Value = "Adass Israel School Elsternwick"
$ws.Cells.Item(25, 2).Value = 17
$ws.Cells.Item(26, 1).Value = "Antonine College Cedar Campus Coburg"
$ws.Cells.Item(26, 2).Value = 10
$ws.Cells.Item(27, 1).Value = "Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh"
$ws.Cells.Item(27, 2).Value = 11
$ws.Cells.Item(28, 1).Value = "Covenant College Bell Post Hill"
$ws.Cells.Item(28, 2).Value = 22
$ws.Cells.Item(29, 1).Value = "Creekside K-9 College Caroline Springs"
$ws.Cells.Item(29, 2).Value = 17
$ws.Cells.Item(30, 1).Value = "Darul Ulum College of Victoria Fawkner October"
$ws.Cells.Item(30, 2).Value = 11
$ws.Cells.Item(31, 1).Value = "Derrimut Primary School Derrimut"
$ws.Cells.Item(31, 2).Value = 11
$ws.Cells.Item(32, 1).Value = "Devon Meadows Primary School Devon Meadows"
$ws.Cells.Item(32, 2).Value = 10
$ws.Cells.Item(33, 1).Value = "Drouin Primary School Drouin"
$ws.Cells.Item(33, 2).Value = 11
$ws.Cells.Item(34, 1).Value = "Exford Primary School Exford"
$ws.Cells.Item(34, 2).Value = 15
$ws.Cells.Item(35, 1).Value = "Flemington Racecourse Flemington"
$ws.Cells.Item(35, 2).Value = 15
$ws.Cells.Item(36, 1).Value = "Gilly's Early Learning Centre Balaclava"
$ws.Cells.Item(36, 2).Value = 10
$ws.Cells.Item(37, 1).Value = "Hazel Glen College Doreen"
$ws.Cells.Item(37, 2).Value = 16
$ws.Cells.Item(38, 1).Value = "Hazelwood North Primary School Hazelwood North"
$ws.Cells.Item(38, 2).Value = 19
$ws.Cells.Item(39, 1).Value = "Ilim College Dallas Main Campus Dallas Oct"
$ws.Cells.Item(39, 2).Value = 17
$ws.Cells.Item(40, 1).Value = "Ilim College Kiewa Campus Dallas"
$ws.Cells.Item(40, 2).Value = 10
$ws.Cells.Item(41, 1).Value = "Islamic College of Melbourne Tarneit Oct Nov"
$ws.Cells.Item(41, 2).Value = 60
$ws.Cells.Item(42, 1).Value = "Lyndhurst Primary School Lyndhurst"
$ws.Cells.Item(42, 2).Value = 11
$ws.Cells.Item(43, 1).Value = "Middle Park Primary School Middle Park"
$ws.Cells.Item(43, 2).Value = 14
$ws.Cells.Item(44, 1).Value = "Morwell Park Primary School Morwell"
$ws.Cells.Item(44, 2).Value = 76
$ws.Cells.Item(45, 1).Value = "Nio Early Learning Adventures Preston"
$ws.Cells.Item(45, 2).Value = 19
$ws.Cells.Item(46, 1).Value = "Pentland Primary School Darley"
$ws.Cells.Item(46, 2).Value = 13
$ws.Cells.Item(47, 1).Value = "Rutherglen Motor Inn and Walkabout Motel Rutherglen"
$ws.Cells.Item(47, 2).Value = 16
$ws.Cells.Item(48, 1).Value = "Sirius College Ibrahim Dellal Campus Sunshine"
$ws.Cells.Item(48, 2).Value = 13
$ws.Cells.Item(49, 1).Value = "Sirius College Shepparton Campus Shepparton"
$ws.Cells.Item(49, 2).Value = 16
$ws.Cells.Item(50, 1).Value = "Social Gathering Woodvale 30 Oct"
$ws.Cells.Item(50, 2).Value = 10
$ws.Cells.Item(51, 1).Value = "Society Restaurant Melbourne"
$ws.Cells.Item(51, 2).Value = 36
$ws.Cells.Item(52, 1).Value = "St Ambrose Parish Primary School Woodend"
$ws.Cells.Item(52, 2).Value = 12
$ws.Cells.Item(53, 1).Value = "St Paul's Primary School Sunshine West"
$ws.Cells.Item(53, 2).Value = 15
$ws.Cells.Item(54, 1).Value = "Stevensville Primary School St Albans"
$ws.Cells.Item(54, 2).Value = 10
$ws.Cells.Item(55, 1).Value = "Stockdale Road Primary School Traralgon"
$ws.Cells.Item(55, 2).Value = 29
$ws.Cells.Item(56, 1).Value = "Supreme Caravans Manufacturing Campbellfield"
$ws.Cells.Item(56, 2).Value = 51
$ws.Cells.Item(57, 1).Value = "Templestowe Park Primary School Templestowe"
$ws.Cells.Item(57, 2).Value = 19
$ws.Cells.Item(58, 1).Value = "The Lake Primary School Cabarita"
$ws.Cells.Item(58, 2).Value = 12
$ws.Cells.Item(59, 1).Value = "Top Yard Rooftop Melbourne"
$ws.Cells.Item(59, 2).Value = 14
$ws.Cells.Item(60, 1).Value = "Truganina P-9 College Truganina"
$ws.Cells.Item(60, 2).Value = 14
$ws.Cells.Item(61, 1).Value = "Tucker Road Bentleigh Primary School Bentleigh"
$ws.Cells.Item(61, 2).Value = 10
$ws.Cells.Item(62, 1).Value = "Warragul Regional College Warragul"
$ws.Cells.Item(62, 2).Value = 22
$ws.Cells.Item(63, 1).Value = "Yeshivah College St Kilda East"
$ws.Cells.Item(63, 2).Value = 24
